$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a serial date value of 45205 for every
# data row (rows 2-490). Update it to 45206 for all of them.
$ws.Range("C2:C490").Value = 45206
